# Add data for 2024-06-29
# Updates 2024 (column K) violent-crime counts across the Citywide Totals,
# By Neighborhood summary, and the affected individual neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('K2').Value = 3819
$ws.Range('K3').Value = 3869
$ws.Range('K4').Value = 792
$ws.Range('K5').Value = 268
$ws.Range('K6').Value = 4409
$ws.Range('K7').Value = 13157

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('K6').Value = 100
$ws.Range('K7').Value = 383
$ws.Range('K8').Value = 889
$ws.Range('K11').Value = 266
$ws.Range('K15').Value = 133
$ws.Range('K19').Value = 406
$ws.Range('K20').Value = 295
$ws.Range('K21').Value = 34
$ws.Range('K24').Value = 43
$ws.Range('K25').Value = 57
$ws.Range('K29').Value = 695
$ws.Range('K31').Value = 144
$ws.Range('K33').Value = 544
$ws.Range('K36').Value = 162
$ws.Range('K37').Value = 446
$ws.Range('K40').Value = 33
$ws.Range('K42').Value = 461
$ws.Range('K43').Value = 119
$ws.Range('K49').Value = 75
$ws.Range('K51').Value = 151
$ws.Range('K53').Value = 179
$ws.Range('K54').Value = 249
$ws.Range('K55').Value = 149
$ws.Range('K60').Value = 84
$ws.Range('K64').Value = 80
$ws.Range('K72').Value = 63
$ws.Range('K73').Value = 119
$ws.Range('K76').Value = 189
$ws.Range('K79').Value = 344
$ws.Range('K83').Value = 280
$ws.Range('K84').Value = 96
$ws.Range('K85').Value = 600
$ws.Range('K88').Value = 149
$ws.Range('K89').Value = 182
$ws.Range('K91').Value = 142
$ws.Range('K95').Value = 220
$ws.Range('K96').Value = 154
$ws.Range('K97').Value = 111
$ws.Range('K99').Value = 233
$ws.Range('K101').Value = 13157

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('K6').Value = 70
$ws.Range('K7').Value = 154

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('K2').Value = 142
$ws.Range('K6').Value = 92
$ws.Range('K7').Value = 383

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('K5').Value = 3
$ws.Range('K7').Value = 266

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('K2').Value = 45
$ws.Range('K6').Value = 58
$ws.Range('K7').Value = 182

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('K2').Value = 215
$ws.Range('K3').Value = 199
$ws.Range('K6').Value = 137
$ws.Range('K7').Value = 600

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('K3').Value = 39
$ws.Range('K7').Value = 179

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('K2').Value = 256
$ws.Range('K3').Value = 270
$ws.Range('K6').Value = 293
$ws.Range('K7').Value = 889

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('K3').Value = 93
$ws.Range('K7').Value = 280

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('K3').Value = 206
$ws.Range('K4').Value = 25
$ws.Range('K6').Value = 156
$ws.Range('K7').Value = 544

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('K2').Value = 69
$ws.Range('K3').Value = 79
$ws.Range('K7').Value = 220

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('K2').Value = 121
$ws.Range('K3').Value = 152
$ws.Range('K6').Value = 134
$ws.Range('K7').Value = 446

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('K4').Value = 14
$ws.Range('K7').Value = 233

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range('K6').Value = 52
$ws.Range('K7').Value = 144

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('K3').Value = 40
$ws.Range('K7').Value = 96

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('K3').Value = 14
$ws.Range('K4').Value = 7
$ws.Range('K7').Value = 75

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('K6').Value = 118
$ws.Range('K7').Value = 249

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('K2').Value = 198
$ws.Range('K4').Value = 38
$ws.Range('K6').Value = 199
$ws.Range('K7').Value = 695

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('K2').Value = 22
$ws.Range('K4').Value = 23

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('K2').Value = 132
$ws.Range('K3').Value = 121
$ws.Range('K4').Value = 15
$ws.Range('K7').Value = 406

$ws = $wb.Worksheets.Item('River North')
$ws.Range('K4').Value = 10
$ws.Range('K7').Value = 189

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('K2').Value = 42
$ws.Range('K7').Value = 100

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('K2').Value = 124
$ws.Range('K7').Value = 461

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('K3').Value = 38
$ws.Range('K7').Value = 149

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('K6').Value = 13
$ws.Range('K7').Value = 43

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('K5').Value = 3
$ws.Range('K6').Value = 33
$ws.Range('K7').Value = 142

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('K6').Value = 18
$ws.Range('K7').Value = 34

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('K2').Value = 118
$ws.Range('K3').Value = 112
$ws.Range('K6').Value = 81
$ws.Range('K7').Value = 344

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('K6').Value = 31
$ws.Range('K7').Value = 80

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('K6').Value = 92
$ws.Range('K7').Value = 295

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('K6').Value = 35
$ws.Range('K7').Value = 162

$ws = $wb.Worksheets.Item('East Side')
$ws.Range('K4').Value = 5
$ws.Range('K7').Value = 57

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range('K4').Value = 9
$ws.Range('K6').Value = 42
$ws.Range('K7').Value = 133

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('K6').Value = 48
$ws.Range('K7').Value = 119

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('K6').Value = 70
$ws.Range('K7').Value = 111

$ws = $wb.Worksheets.Item('United Center')
$ws.Range('K3').Value = 41
$ws.Range('K6').Value = 69
$ws.Range('K7').Value = 149

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('K3').Value = 43
$ws.Range('K7').Value = 151

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('K2').Value = 28
$ws.Range('K7').Value = 84

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('K3').Value = 33
$ws.Range('K7').Value = 119

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('K6').Value = 34
$ws.Range('K7').Value = 63

$ws = $wb.Worksheets.Item('Hegewisch')
$ws.Range('K2').Value = 12
$ws.Range('K7').Value = 33
